$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '43.030.73'
$ws.Range("D3").Value = '2.544.44'
$ws.Range("E3").Value = '  +0.77%  '
$ws.Range("D4").Value = '''0.999'
$ws.Range("E4").Value = '  -0.12%  '
$ws.Range("D5").Value = '''318.84'
$ws.Range("E5").Value = '  +0.86%  '
$ws.Range("D6").Value = '''97.86'
$ws.Range("E6").Value = '  +3.18%  '
$ws.Range("D7").Value = '''0.576'
$ws.Range("E7").Value = '  -0.66%  '
$ws.Range("E8").Value = '  +0.01%  '
$ws.Range("D9").Value = '''0.537'
$ws.Range("E9").Value = '  -0.15%  '
$ws.Range("D10").Value = '''36.39'
$ws.Range("E10").Value = '  +0.43%  '
$ws.Range("D11").Value = '''0.0820'
$ws.Range("E11").Value = '  +1.00%  '
$ws.Range("D12").Value = '''7.66'
$ws.Range("E12").Value = '  +1.11%  '
$ws.Range("E13").Value = '  -3.09%  '
$ws.Range("D14").Value = '2.935.50'
$ws.Range("E14").Value = '  +0.66%  '
$ws.Range("D15").Value = '2.560.03'
$ws.Range("E15").Value = '  +0.66%  '
$ws.Range("D16").Value = '''15.22'
$ws.Range("E16").Value = '  -2.13%  '
$ws.Range("D17").Value = '''0.853'
$ws.Range("E17").Value = '  -1.36%  '
$ws.Range("D18").Value = '43.075.71'
$ws.Range("E18").Value = '  +0.71%  '
$ws.Range("D19").Value = '''6.88'
$ws.Range("E19").Value = '  +4.81%  '
$ws.Range("D20").Value = '''12.88'
$ws.Range("E20").Value = '  -1.15%  '
$ws.Range("D21").Value = '0.0₃0970'
$ws.Range("E21").Value = '  +0.11%  '
$ws.Range("D22").Value = '''70.04'
$ws.Range("E22").Value = '  -1.83%  '
$ws.Range("D23").Value = '''255.09'
$ws.Range("E23").Value = '  -0.13%  '
$ws.Range("D24").Value = '''2.96'
$ws.Range("E24").Value = '  -0.83%  '
$ws.Range("D25").Value = '''2.06'
$ws.Range("E25").Value = '  +0.86%  '
$ws.Range("D26").Value = '''26.68'
$ws.Range("E26").Value = '  -3.98%  '
$ws.Range("E27").Value = '  +0.61%  '
$ws.Range("E28").Value = '  +5.32%  '
$ws.Range("D29").Value = '''41.11'
$ws.Range("E29").Value = '  +4.59%  '
$ws.Range("D30").Value = '''10.49'
$ws.Range("E30").Value = '  +4.30%  '
$ws.Range("D31").Value = '''5.95'
$ws.Range("E31").Value = '  +0.71%  '
$ws.Range("D32").Value = '''158.42'
$ws.Range("E32").Value = '  +1.36%  '
$ws.Range("D33").Value = '''2.18'
$ws.Range("E33").Value = '  +3.48%  '
$ws.Range("E34").Value = '  +0.96%  '
$ws.Range("E35").Value = '  +4.52%  '
$ws.Range("D36").Value = '''19.06'
$ws.Range("E36").Value = '  -4.13%  '
$ws.Range("E37").Value = '  +1.19%  '
$ws.Range("E38").Value = '  +0.25%  '
$ws.Range("E39").Value = '  +14.15%  '
$ws.Range("E40").Value = '  -0.55%  '
$ws.Range("D41").Value = '''22.04'
$ws.Range("E41").Value = '  -10.86%  '
$ws.Range("D42").Value = '''3.86'
$ws.Range("E42").Value = '  +0.19%  '
$ws.Range("E43").Value = '  +0.41%  '
$ws.Range("E44").Value = '  +0.29%  '
$ws.Range("E45").Value = '  -2.71%  '
$ws.Range("D46").Value = '2.026.82'
$ws.Range("E46").Value = '  -1.25%  '
$ws.Range("D47").Value = '''9.14'
$ws.Range("E47").Value = '  +3.47%  '
$ws.Range("D48").Value = '''84.80'
$ws.Range("E48").Value = '  -2.02%  '
$ws.Range("D49").Value = '''76.37'
$ws.Range("E49").Value = '  +2.79%  '
$ws.Range("D50").Value = '''107.03'
$ws.Range("E50").Value = '  +4.83%  '
$ws.Range("D51").Value = '2.789.70'
$ws.Range("E51").Value = '  +0.75%  '

# The apostrophe-prefixed assignments above force Excel to keep those
# numeric-looking values as text, but they also pick up a "quote prefix"
# cell style. Paste formats from an unstyled cell (D13, never re-styled)
# back onto them so the workbook-level style table is left untouched.
$ws.Range("D13").Copy() | Out-Null
$ws.Range("D4").PasteSpecial(-4122)
$ws.Range("D5").PasteSpecial(-4122)
$ws.Range("D6").PasteSpecial(-4122)
$ws.Range("D7").PasteSpecial(-4122)
$ws.Range("D9").PasteSpecial(-4122)
$ws.Range("D10").PasteSpecial(-4122)
$ws.Range("D11").PasteSpecial(-4122)
$ws.Range("D12").PasteSpecial(-4122)
$ws.Range("D16").PasteSpecial(-4122)
$ws.Range("D17").PasteSpecial(-4122)
$ws.Range("D19").PasteSpecial(-4122)
$ws.Range("D20").PasteSpecial(-4122)
$ws.Range("D22").PasteSpecial(-4122)
$ws.Range("D23").PasteSpecial(-4122)
$ws.Range("D24").PasteSpecial(-4122)
$ws.Range("D25").PasteSpecial(-4122)
$ws.Range("D26").PasteSpecial(-4122)
$ws.Range("D29").PasteSpecial(-4122)
$ws.Range("D30").PasteSpecial(-4122)
$ws.Range("D31").PasteSpecial(-4122)
$ws.Range("D32").PasteSpecial(-4122)
$ws.Range("D33").PasteSpecial(-4122)
$ws.Range("D36").PasteSpecial(-4122)
$ws.Range("D41").PasteSpecial(-4122)
$ws.Range("D42").PasteSpecial(-4122)
$ws.Range("D47").PasteSpecial(-4122)
$ws.Range("D48").PasteSpecial(-4122)
$ws.Range("D49").PasteSpecial(-4122)
$ws.Range("D50").PasteSpecial(-4122)
$excel.CutCopyMode = $false
